$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '30.284.13'
$ws.Range('E2').NumberFormat = "@"
$ws.Range('E2').Value = '  -0.04%  '

$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '1.866.45'
$ws.Range('E3').NumberFormat = "@"
$ws.Range('E3').Value = '  -0.81%  '

$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.0000'
$ws.Range('E4').NumberFormat = "@"
$ws.Range('E4').Value = '  -0.06%  '

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '235.27'
$ws.Range('E5').NumberFormat = "@"
$ws.Range('E5').Value = '  -0.34%  '

$ws.Range('E6').NumberFormat = "@"
$ws.Range('E6').Value = '  -0.02%  '

$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.4679'
$ws.Range('E7').NumberFormat = "@"
$ws.Range('E7').Value = '  -0.13%  '

$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.2839'
$ws.Range('E8').NumberFormat = "@"
$ws.Range('E8').Value = '  +0.42%  '

$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.06522'
$ws.Range('E9').NumberFormat = "@"
$ws.Range('E9').Value = '  -0.97%  '

$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '21.29'
$ws.Range('E10').NumberFormat = "@"
$ws.Range('E10').Value = '  +3.10%  '

$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.07868'
$ws.Range('E11').NumberFormat = "@"
$ws.Range('E11').Value = '  +1.28%  '

$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '97.39'
$ws.Range('E12').NumberFormat = "@"
$ws.Range('E12').Value = '  -0.35%  '

$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '1.880.70'
$ws.Range('E13').NumberFormat = "@"
$ws.Range('E13').Value = '  -0.16%  '

$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '5.095'
$ws.Range('E14').NumberFormat = "@"
$ws.Range('E14').Value = '  +0.40%  '

$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.6743'
$ws.Range('E15').NumberFormat = "@"
$ws.Range('E15').Value = '  +0.02%  '

$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '279.68'
$ws.Range('E16').NumberFormat = "@"
$ws.Range('E16').Value = '  -1.29%  '

$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '30.279.89'
$ws.Range('E17').NumberFormat = "@"
$ws.Range('E17').Value = '  -0.08%  '

$ws.Range('E18').NumberFormat = "@"
$ws.Range('E18').Value = '  +0.09%  '

$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '5.495'
$ws.Range('E19').NumberFormat = "@"
$ws.Range('E19').Value = '  +1.61%  '

$ws.Range('E20').NumberFormat = "@"
$ws.Range('E20').Value = '  +0.48%  '

$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '2.124.20'

$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '0.000007279'
$ws.Range('E22').NumberFormat = "@"
$ws.Range('E22').Value = '  +0.07%  '

$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '1.0000'
$ws.Range('E23').NumberFormat = "@"
$ws.Range('E23').Value = '  -0.16%  '

$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '6.157'
$ws.Range('E24').NumberFormat = "@"
$ws.Range('E24').Value = '  -0.38%  '

$ws.Range('B25').Value = 'Cosmos'
$ws.Range('C25').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '9.188'
$ws.Range('E25').NumberFormat = "@"
$ws.Range('E25').Value = '  -2.09%  '

$ws.Range('B26').Value = 'Monero'
$ws.Range('C26').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '165.22'
$ws.Range('E26').NumberFormat = "@"
$ws.Range('E26').Value = '  -1.66%  '

$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '19.11'
$ws.Range('E27').NumberFormat = "@"
$ws.Range('E27').Value = '  -0.39%  '

$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '1.927'
$ws.Range('E28').NumberFormat = "@"
$ws.Range('E28').Value = '  -3.17%  '

$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '1.376'
$ws.Range('E29').NumberFormat = "@"
$ws.Range('E29').Value = '  +0.29%  '

$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '0.09655'
$ws.Range('E30').NumberFormat = "@"
$ws.Range('E30').Value = '  -0.16%  '

$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '4.383'
$ws.Range('E31').NumberFormat = "@"
$ws.Range('E31').Value = '  +0.17%  '

$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '1.476'
$ws.Range('E32').NumberFormat = "@"
$ws.Range('E32').Value = '  -0.15%  '

$ws.Range('E33').NumberFormat = "@"
$ws.Range('E33').Value = '  -0.78%  '

$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.04707'
$ws.Range('E34').NumberFormat = "@"
$ws.Range('E34').Value = '  +0.81%  '

$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '1.126'
$ws.Range('E35').NumberFormat = "@"
$ws.Range('E35').Value = '  +2.59%  '

$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.7058'
$ws.Range('E36').NumberFormat = "@"
$ws.Range('E36').Value = '  +0.16%  '

$ws.Range('E37').NumberFormat = "@"
$ws.Range('E37').Value = '  +0.47%  '

$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.01855'
$ws.Range('E38').NumberFormat = "@"
$ws.Range('E38').Value = '  -0.79%  '

$ws.Range('B39').Value = 'FraxShare'
$ws.Range('C39').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '6.250'
$ws.Range('E39').NumberFormat = "@"
$ws.Range('E39').Value = '  -5.23%  '

$ws.Range('B40').Value = 'MXToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '2.528'
$ws.Range('E40').NumberFormat = "@"
$ws.Range('E40').Value = '  +0.14%  '

$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '73.52'
$ws.Range('E41').NumberFormat = "@"
$ws.Range('E41').Value = '  +2.14%  '

$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '1.946'
$ws.Range('E42').NumberFormat = "@"
$ws.Range('E42').Value = '  -1.11%  '

$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.8473'
$ws.Range('E43').NumberFormat = "@"
$ws.Range('E43').Value = '  -2.15%  '

$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.4174'
$ws.Range('E44').NumberFormat = "@"
$ws.Range('E44').Value = '  -0.05%  '

$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.9998'

$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '103.82'
$ws.Range('E46').NumberFormat = "@"
$ws.Range('E46').Value = '  +0.79%  '

$ws.Range('E47').NumberFormat = "@"
$ws.Range('E47').Value = '  -1.26%  '

$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '9.164'
$ws.Range('E48').NumberFormat = "@"
$ws.Range('E48').Value = '  -0.26%  '

$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '935.10'
$ws.Range('E49').NumberFormat = "@"
$ws.Range('E49').Value = '  -4.84%  '

$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '34.06'
$ws.Range('E50').NumberFormat = "@"
$ws.Range('E50').Value = '  +0.51%  '

$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.1123'
$ws.Range('E51').NumberFormat = "@"
$ws.Range('E51').Value = '  -1.80%  '
